# "create new packaging method" -- add a new PackagedProductDefinition row
# describing the 30-tablet (1 bottle of 30) bottle packaging, alongside the
# existing 90-tablet (3 bottles of 30) packaging on row 2, and leave the UI
# focused on the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PackagedProductDefinition")

# Duplicate the existing packaging row (B2:N2) down into row 3, carrying over
# styles/number formats, then patch the two cells that describe the new
# packaging method (name + quantity text).
$ws.Range("B2:N2").Copy($ws.Range("B3"))
$ws.Range("H3").Value2 = "30 tablets"
$ws.Range("C3").Value2 = "Biktarvy 50 mg/200 mg/25 mg film-coated tablets 30 (1 bottle of 30) film-coated tablets"

# Move the active selection/tab to the PackagedProductDefinition sheet, on
# the newly-created row.
$ws.Activate()
$ws.Range("C11").Select()
